# Generate Report for Handoff
# b.md has moved from "Handed back: in sync with en-US" to "Ready for handoff":
# a new handoff package (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf) was
# cut for it on both locales, so the Overview roll-up and each locale sheet's
# row for b.md need their Status / Latest Handoff File / Latest Handoff
# Datetime columns refreshed to reflect the new handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is b.md ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-26-12 18:26:55"

# --- zh-cn sheet: row 3 is b.md ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-12 18:26:51"
$zh.Hyperlinks.Item(8).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# --- de-de sheet: row 3 is b.md ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("E3").Value = "2016-03-12 18:26:55"
$de.Hyperlinks.Item(8).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
